$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 becomes the old S2 values (name changes S1 -> S2), with D5 now 0 instead of 3
$ws.Range("A5").Value = "S2"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 22.92
$ws.Range("D5").Value = 0

# Row 6 becomes the old S1 values (name changes S2 -> S1), with updated E6/F6
$ws.Range("A6").Value = "S1"
$ws.Range("B6").Value = 19.64
$ws.Range("C6").Value = 60.01
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 43.64333333333333
$ws.Range("F6").Value = 1.38
